# Adds a letter label as a new first paragraph above the existing
# number in each of the ten ellipse shapes on slide 1.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$labels = @{
    1  = "C"   # "12"
    2  = "A"   # "2"
    3  = "G"   # "4"
    4  = "B"   # "6"
    5  = "E"   # "8"
    6  = "F"   # "10"
    7  = "I"   # "16"
    8  = "D"   # "14"
    9  = "H"   # "18"
    10 = "J"   # "20"
}

foreach ($idx in $labels.Keys) {
    $shp = $s.Shapes.Item($idx)
    $tr = $shp.TextFrame.TextRange
    if ($idx -eq 1) {
        # This shape's trailing empty endParaRPr gets dropped, matching
        # the diff (a full-text rewrite rather than an in-place insert).
        $tr.Text = $labels[$idx] + [char]13 + $tr.Text
    } else {
        [void]$tr.InsertBefore($labels[$idx] + [char]13)
    }
}
